# Nuevo formato 15 jun 2021
# Update the "row 10" statistics block on each of the three parcial sheets.

$wb = $excel.ActiveWorkbook

# --- 1er Parcial ---
$ws1 = $wb.Worksheets.Item("1er Parcial")
$ws1.Range("E10").Value = 22
$ws1.Range("F10").Value = 11
$ws1.Range("G10").Value = 66.67
$ws1.Range("H10").Value = 33.33
$ws1.Range("J10").Value = 0
$ws1.Range("K10").Value = 0

# --- 2o Parcial ---
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Range("E10").Value = 30
$ws2.Range("F10").Value = 3
$ws2.Range("G10").Value = 90.91
$ws2.Range("H10").Value = 9.09
$ws2.Range("I10").Value = 7.6
$ws2.Range("J10").Value = 0
$ws2.Range("K10").Value = 0

# --- 3er Parcial ---
$ws3 = $wb.Worksheets.Item("3er Parcial")
$ws3.Range("E10").Value = 30
$ws3.Range("F10").Value = 3
$ws3.Range("G10").Value = 90.91
$ws3.Range("H10").Value = 9.09
$ws3.Range("I10").Value = 7.5
$ws3.Range("J10").Value = 0
$ws3.Range("K10").Value = 0
